{"js": "// The document's sole paragraph reads \"Version 2.\" and must become\n// \"Version 1.\" (per the diff: \"Versi\"+\"on\" merge into a single \"Version\"\n// run, and \" 2\"+\".\" merge into a single \" 1.\" run that now sits entirely\n// before the _GoBack bookmark).\nconst body = context.document.body;\n\n// Step 1: re-write \"Version\" (spans the original \"Versi\"/\"on\" runs) so the\n// two runs collapse into a single run with the same text.\nconst versionHits = body.search(\"Version\", { matchCase: true });\nversionHits.load(\"items\");\nawait context.sync();\n\nif (versionHits.items.length > 0) {\n  versionHits.items[0].insertText(\"Version\", \"Replace\");\n}\nawait context.sync();\n\n// Step 2: re-write \" 2.\" (spans the original \" 2\" run, the _GoBack\n// bookmark, and the trailing \".\" run) as \" 1.\", bumping the version\n// number and collapsing the two text runs into one that now precedes the\n// bookmark.\nconst versionNumberHits = body.search(\" 2.\", { matchCase: true });\nversionNumberHits.load(\"items\");\nawait context.sync();\n\nif (versionNumberHits.items.length > 0) {\n  versionNumberHits.items[0].insertText(\" 1.\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document's sole paragraph reads \"Version 2.\" and must become\n# \"Version 1.\" (per the diff: \"Versi\"+\"on\" merge into a single \"Version\"\n# run, and \" 2\"+\".\" merge into a single \" 1.\" run that now sits entirely\n# before the _GoBack bookmark).\n\n# Step 1: re-write \"Version\" (spans the original \"Versi\"/\"on\" runs) so the\n# two runs collapse into a single run with the same text.\n$find1 = $d.Content.Find\n$find1.Text = \"Version\"\n$find1.Replacement.Text = \"Version\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# Step 2: bump the version number, \"2\" -> \"1\". This stays inside the\n# \" 2\" run (it doesn't reach the _GoBack bookmark), so the bookmark is\n# left untouched.\n$find2 = $d.Content.Find\n$find2.Text = \"2\"\n$find2.Replacement.Text = \"1\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n\n# Step 3: the trailing \".\" currently sits in its own run right after the\n# _GoBack bookmark. Delete it there and re-insert it immediately before\n# the bookmark so it rejoins the \" 1\" run as a single \" 1.\" run, matching\n# the target structure and keeping the bookmark in place.\n$bm = $d.Bookmarks(\"_GoBack\")\n$bmStart = $bm.Range.Start\n\n$period = $d.Range($bmStart, $bmStart + 1)\nif ($period.Text -eq \".\") {\n    $period.Delete()\n    $gap = $d.Range($bmStart, $bmStart)\n    $gap.InsertBefore(\".\")\n}\n"}
